$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.634.25"
$ws.Range("E2").Value = "  +3.94%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.422.79"
$ws.Range("E3").Value = "  +2.61%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.22"
$ws.Range("E5").Value = "  +4.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.23"
$ws.Range("E6").Value = "  +6.03%  "

# Row 7
$ws.Range("E7").Value = "  +2.48%  "

# Row 8
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("E9").Value = "  +8.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.46"
$ws.Range("E10").Value = "  +4.01%  "

# Row 11
$ws.Range("E11").Value = "  +2.05%  "

# Row 12
$ws.Range("E12").Value = "  +2.99%  "

# Row 13
$ws.Range("E13").Value = "  -2.01%  "

# Row 14
$ws.Range("E14").Value = "  +3.86%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.800.82"
$ws.Range("E15").Value = "  +2.71%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.420.43"
$ws.Range("E16").Value = "  +2.97%  "

# Row 17
$ws.Range("E17").Value = "  +5.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.483.02"
$ws.Range("E18").Value = "  +3.67%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.46"
$ws.Range("E19").Value = "  +4.59%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.41"
$ws.Range("E20").Value = "  +2.53%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0923"
$ws.Range("E21").Value = "  +4.36%  "

# Row 22
$ws.Range("E22").Value = "  +1.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.72"
$ws.Range("E23").Value = "  +3.29%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.29"
$ws.Range("E24").Value = "  +5.86%  "

# Row 25
$ws.Range("E25").Value = "  +2.01%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.27"
$ws.Range("E27").Value = "  +3.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("E28").Value = "  -3.54%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.53"
$ws.Range("E29").Value = "  +2.35%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.35"
$ws.Range("E30").Value = "  +4.18%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.47"
$ws.Range("E31").Value = "  +1.43%  "

# Row 32
$ws.Range("E32").Value = "  +20.49%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.41"
$ws.Range("E33").Value = "  +10.84%  "

# Row 34
$ws.Range("E34").Value = "  +3.79%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0778"
$ws.Range("E35").Value = "  +8.84%  "

# Row 36
$ws.Range("E36").Value = "  +0.19%  "

# Row 37
$ws.Range("E37").Value = "  +3.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.51"
$ws.Range("E38").Value = "  +4.95%  "

# Row 39
$ws.Range("E39").Value = "  +1.95%  "

# Row 40
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.109"
$ws.Range("E40").Value = "  +1.89%  "

# Row 41
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "119.95"
$ws.Range("E41").Value = "  -5.10%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  -2.82%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.02"
$ws.Range("E43").Value = "  -1.18%  "

# Row 44
$ws.Range("E44").Value = "  +4.74%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.940.27"
$ws.Range("E45").Value = "  +0.75%  "

# Row 46
$ws.Range("E46").Value = "  +1.56%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.94"
$ws.Range("E47").Value = "  +8.99%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.48"
$ws.Range("E48").Value = "  +3.04%  "

# Row 49
$ws.Range("E49").Value = "  +11.16%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.82"
$ws.Range("E50").Value = "  +6.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.57"
$ws.Range("E51").Value = "  +6.70%  "
